$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5566744
$ws.Range("J43").Value = 19254.555
$ws.Range("L43").Value = 19254.555
$ws.Range("N43").Value = -19392.555
$ws.Range("H74").Value = 12842.833
$ws.Range("I74").Value = 11411.4
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 11411.4
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -10475.4
$ws.Range("N74").Value = -21872
$ws.Range("H77").Value = 12842.833
$ws.Range("I77").Value = 11411.4
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 57057
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -52377
$ws.Range("N77").Value = -109360
$ws.Range("H80").Value = 647.5
$ws.Range("I80").Value = 686.6
$ws.Range("J80").Value = 582.3333
$ws.Range("K80").Value = 2059.8
$ws.Range("L80").Value = 1746.9999
$ws.Range("M80").Value = -1061.8
$ws.Range("N80").Value = -3742.9999
$ws.Range("H83").Value = 647.5
$ws.Range("I83").Value = 686.6
$ws.Range("J83").Value = 582.3333
$ws.Range("K83").Value = 6179.400000000001
$ws.Range("L83").Value = 5240.9997
$ws.Range("M83").Value = -1187.400000000001
$ws.Range("N83").Value = -15224.9997
$ws.Range("H98").Value = 1999
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H111").Value = 2133
$ws.Range("I111").Value = 1250
$ws.Range("J111").Value = 3016
$ws.Range("K111").Value = 3750
$ws.Range("L111").Value = 9048
$ws.Range("M111").Value = -683
$ws.Range("N111").Value = -15182
$ws.Range("H116").Value = 4999.25
$ws.Range("I116").Value = 5165.6665
$ws.Range("K116").Value = 5165.6665
$ws.Range("M116").Value = -1723.6665
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H125").Value = 1798.3334
$ws.Range("I125").Value = 1847.5
$ws.Range("K125").Value = 16627.5
$ws.Range("M125").Value = -14167.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 940.95
$ws.Range("J2").Value = 1052.4445
$ws.Range("L2").Value = 1052.4445
$ws.Range("N2").Value = -1278.4445
$ws.Range("H45").Value = 2721.4348
$ws.Range("I45").Value = 2071.2307
$ws.Range("K45").Value = 2071.2307
$ws.Range("M45").Value = -1694.2307
$ws.Range("H63").Value = 8985
$ws.Range("I63").Value = 1966
$ws.Range("K63").Value = 1966
$ws.Range("M63").Value = -1280
$ws.Range("H66").Value = 8985
$ws.Range("I66").Value = 1966
$ws.Range("K66").Value = 9830
$ws.Range("M66").Value = -6398
$ws.Range("H97").Value = 1929.1428
$ws.Range("I97").Value = 1917.3334
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1917.3334
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1421.3334
$ws.Range("N97").Value = -2992
$ws.Range("H116").Value = 940.95
$ws.Range("J116").Value = 1052.4445
$ws.Range("L116").Value = 1052.4445
$ws.Range("N116").Value = -5640.4445
$ws.Range("H122").Value = 1827.48
$ws.Range("I122").Value = 1402.95
$ws.Range("K122").Value = 4208.85
$ws.Range("M122").Value = -1758.85

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 940.95
$ws.Range("J3").Value = 1052.4445
$ws.Range("L3").Value = 1052.4445
$ws.Range("N3").Value = -1280.4445
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H26").Value = 37500
$ws.Range("I26").Value = 37500
$ws.Range("K26").Value = 37500
$ws.Range("M26").Value = -37208
$ws.Range("H54").Value = 3717.5
$ws.Range("I54").Value = 3717.5
$ws.Range("K54").Value = 3717.5
$ws.Range("M54").Value = -3233.5
$ws.Range("H64").Value = 1010.5
$ws.Range("I64").Value = 980.6667
$ws.Range("K64").Value = 980.6667
$ws.Range("M64").Value = -755.6667
$ws.Range("H67").Value = 1010.5
$ws.Range("I67").Value = 980.6667
$ws.Range("K67").Value = 980.6667
$ws.Range("M67").Value = -200.6667
$ws.Range("H135").Value = 182882.5
$ws.Range("J135").Value = 182882.5
$ws.Range("L135").Value = 182882.5
$ws.Range("N135").Value = -193022.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7306.96
$ws.Range("J31").Value = 8043.5713
$ws.Range("L31").Value = 8043.5713
$ws.Range("N31").Value = -8633.5713
$ws.Range("H34").Value = 7306.96
$ws.Range("J34").Value = 8043.5713
$ws.Range("L34").Value = 8043.5713
$ws.Range("N34").Value = -8447.5713
$ws.Range("H98").Value = 38780
$ws.Range("J98").Value = 38780
$ws.Range("L98").Value = 38780
$ws.Range("N98").Value = -43272
$ws.Range("H132").Value = 3530
$ws.Range("I132").Value = 3530
$ws.Range("K132").Value = 10590
$ws.Range("M132").Value = -8060

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 989
$ws.Range("I3").Value = 989
$ws.Range("K3").Value = 2967
$ws.Range("M3").Value = -2855
$ws.Range("H5").Value = 668.7
$ws.Range("I5").Value = 549
$ws.Range("K5").Value = 1647
$ws.Range("M5").Value = -1535
$ws.Range("H7").Value = 274.91666
$ws.Range("I7").Value = 303.16666
$ws.Range("J7").Value = 246.66667
$ws.Range("K7").Value = 909.4999799999999
$ws.Range("L7").Value = 740.00001
$ws.Range("M7").Value = -797.4999799999999
$ws.Range("N7").Value = -964.00001
$ws.Range("H68").Value = 1951.6
$ws.Range("I68").Value = 1514.25
$ws.Range("J68").Value = 2243.1667
$ws.Range("K68").Value = 4542.75
$ws.Range("L68").Value = 6729.500100000001
$ws.Range("M68").Value = -3731.75
$ws.Range("N68").Value = -8351.500100000001
$ws.Range("H71").Value = 1951.6
$ws.Range("I71").Value = 1514.25
$ws.Range("J71").Value = 2243.1667
$ws.Range("K71").Value = 13628.25
$ws.Range("L71").Value = 20188.5003
$ws.Range("M71").Value = -9572.25
$ws.Range("N71").Value = -28300.5003
$ws.Range("H86").Value = 668.2727
$ws.Range("J86").Value = 721.5714
$ws.Range("L86").Value = 2164.7142
$ws.Range("N86").Value = -4536.7142
$ws.Range("H89").Value = 668.2727
$ws.Range("J89").Value = 721.5714
$ws.Range("L89").Value = 6494.1426
$ws.Range("N89").Value = -18350.1426
$ws.Range("H135").Value = 668.7
$ws.Range("I135").Value = 549
$ws.Range("K135").Value = 4941
$ws.Range("M135").Value = -2406

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 97.85714
$ws.Range("I2").Value = 52.666668
$ws.Range("K2").Value = 52.666668
$ws.Range("M2").Value = 60.333332
$ws.Range("H97").Value = 900
$ws.Range("I97").Value = 800
$ws.Range("K97").Value = 800
$ws.Range("M97").Value = -304
$ws.Range("H102").Value = 1123.5625
$ws.Range("I102").Value = 1123.5625
$ws.Range("K102").Value = 1123.5625
$ws.Range("M102").Value = 498.4375
$ws.Range("H132").Value = 4091.2727
$ws.Range("I132").Value = 3964.1667
$ws.Range("K132").Value = 11892.5001
$ws.Range("M132").Value = -9362.500100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1465.0526
$ws.Range("I55").Value = 933.8
$ws.Range("K55").Value = 933.8
$ws.Range("M55").Value = -760.8
$ws.Range("H68").Value = 8888.444
$ws.Range("I68").Value = 8599.4
$ws.Range("K68").Value = 8599.4
$ws.Range("M68").Value = -7850.4
$ws.Range("H71").Value = 8888.444
$ws.Range("I71").Value = 8599.4
$ws.Range("K71").Value = 42997
$ws.Range("M71").Value = -39253
$ws.Range("H82").Value = 7250
$ws.Range("I82").Value = 5750
$ws.Range("K82").Value = 5750
$ws.Range("M82").Value = -5389
$ws.Range("H85").Value = 7250
$ws.Range("I85").Value = 5750
$ws.Range("K85").Value = 5750
$ws.Range("M85").Value = -4502
$ws.Range("H93").Value = 891.6667
$ws.Range("I93").Value = 837.5
$ws.Range("K93").Value = 837.5
$ws.Range("M93").Value = 410.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 5000.75
$ws.Range("I132").Value = 5001
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 15003
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -12473
$ws.Range("N132").Value = -20060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 34570
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H122").Value = 4517.857
$ws.Range("I122").Value = 4480.769
$ws.Range("K122").Value = 13442.307
$ws.Range("M122").Value = -10992.307
$ws.Range("H126").Value = 5395.8237
$ws.Range("I126").Value = 3203.3333
$ws.Range("K126").Value = 9609.999899999999
$ws.Range("M126").Value = -7139.999899999999
